$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 200899.4
$ws.Range("I62").Value = 250874.75
$ws.Range("J62").Value = 998
$ws.Range("K62").Value = 250874.75
$ws.Range("L62").Value = 998
$ws.Range("M62").Value = -250250.75
$ws.Range("N62").Value = -2246
$ws.Range("H64").Value = 5812
$ws.Range("I64").Value = 4499.5
$ws.Range("J64").Value = 6249.5
$ws.Range("K64").Value = 4499.5
$ws.Range("L64").Value = 6249.5
$ws.Range("M64").Value = -4251.5
$ws.Range("N64").Value = -6745.5
$ws.Range("H65").Value = 200899.4
$ws.Range("I65").Value = 250874.75
$ws.Range("J65").Value = 998
$ws.Range("K65").Value = 1254373.75
$ws.Range("L65").Value = 4990
$ws.Range("M65").Value = -1251253.75
$ws.Range("N65").Value = -11230
$ws.Range("H67").Value = 5812
$ws.Range("I67").Value = 4499.5
$ws.Range("J67").Value = 6249.5
$ws.Range("K67").Value = 4499.5
$ws.Range("L67").Value = 6249.5
$ws.Range("M67").Value = -3641.5
$ws.Range("N67").Value = -7965.5
$ws.Range("H97").Value = 3102.5
$ws.Range("J97").Value = 3102.5
$ws.Range("L97").Value = 9307.5
$ws.Range("N97").Value = -10299.5
$ws.Range("H98").Value = 1576.8846
$ws.Range("I98").Value = 1625
$ws.Range("K98").Value = 1625
$ws.Range("M98").Value = -127
$ws.Range("H100").Value = 6620.3784
$ws.Range("J100").Value = 7630.433
$ws.Range("L100").Value = 7630.433
$ws.Range("N100").Value = -8712.433000000001
$ws.Range("H112").Value = 1611.5
$ws.Range("J112").Value = 1777.8462
$ws.Range("L112").Value = 5333.5386
$ws.Range("N112").Value = -7549.5386
$ws.Range("H113").Value = 6614.385
$ws.Range("J113").Value = 7862.25
$ws.Range("L113").Value = 7862.25
$ws.Range("N113").Value = -14370.25
$ws.Range("H115").Value = 1912.5
$ws.Range("I115").Value = 400
$ws.Range("J115").Value = 3425
$ws.Range("K115").Value = 1200
$ws.Range("L115").Value = 10275
$ws.Range("M115").Value = 367
$ws.Range("N115").Value = -13409
$ws.Range("H118").Value = 1510.6666
$ws.Range("I118").Value = 1212.8
$ws.Range("K118").Value = 3638.4
$ws.Range("M118").Value = -1981.4
$ws.Range("H122").Value = 1576.8846
$ws.Range("I122").Value = 1625
$ws.Range("K122").Value = 4875
$ws.Range("M122").Value = -2425
$ws.Range("H127").Value = 1599.5
$ws.Range("I127").Value = 759.4
$ws.Range("J127").Value = 2999.6667
$ws.Range("K127").Value = 2278.2
$ws.Range("L127").Value = 8999.000100000001
$ws.Range("M127").Value = 2681.8
$ws.Range("N127").Value = -18919.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 100000
$ws.Range("I34").Value = 100000
$ws.Range("K34").Value = 100000
$ws.Range("M34").Value = -99729
$ws.Range("H97").Value = 5155.6665
$ws.Range("I97").Value = 1438.1666
$ws.Range("K97").Value = 1438.1666
$ws.Range("M97").Value = -942.1666
$ws.Range("H110").Value = 1338.7368
$ws.Range("I110").Value = 1211.1333
$ws.Range("J110").Value = 1817.25
$ws.Range("K110").Value = 1211.1333
$ws.Range("L110").Value = 1817.25
$ws.Range("M110").Value = 833.8667
$ws.Range("N110").Value = -5907.25
$ws.Range("H122").Value = 2614.875
$ws.Range("I122").Value = 2504.8076
$ws.Range("K122").Value = 7514.4228
$ws.Range("M122").Value = -5064.4228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 85000
$ws.Range("J51").Value = 85000
$ws.Range("L51").Value = 85000
$ws.Range("N51").Value = -85982
$ws.Range("H107").Value = 2842.3635
$ws.Range("I107").Value = 2909.95
$ws.Range("K107").Value = 2909.95
$ws.Range("M107").Value = -989.9499999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 203.38461
$ws.Range("I7").Value = 197.125
$ws.Range("J7").Value = 213.4
$ws.Range("K7").Value = 197.125
$ws.Range("L7").Value = 213.4
$ws.Range("M7").Value = -84.125
$ws.Range("N7").Value = -439.4
$ws.Range("H16").Value = 678.8333
$ws.Range("I16").Value = 678.8333
$ws.Range("K16").Value = 678.8333
$ws.Range("M16").Value = -391.8333
$ws.Range("H105").Value = 6798.8
$ws.Range("I105").Value = 6798.8
$ws.Range("K105").Value = 6798.8
$ws.Range("M105").Value = -5051.8
$ws.Range("H107").Value = 878.5714
$ws.Range("I107").Value = 885.1667
$ws.Range("J107").Value = 839
$ws.Range("K107").Value = 885.1667
$ws.Range("L107").Value = 839
$ws.Range("M107").Value = 1034.8333
$ws.Range("N107").Value = -4679
$ws.Range("H113").Value = 678.8333
$ws.Range("I113").Value = 678.8333
$ws.Range("K113").Value = 678.8333
$ws.Range("M113").Value = 1491.1667
$ws.Range("H132").Value = 3211
$ws.Range("I132").Value = 3312.375
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 9937.125
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -7407.125
$ws.Range("N132").Value = -12260

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1295.4
$ws.Range("I118").Value = 490.875
$ws.Range("J118").Value = 4513.5
$ws.Range("K118").Value = 1472.625
$ws.Range("L118").Value = 13540.5
$ws.Range("M118").Value = -229.625
$ws.Range("N118").Value = -16026.5
$ws.Range("H131").Value = 2509.25
$ws.Range("J131").Value = 2582
$ws.Range("L131").Value = 7746
$ws.Range("N131").Value = -17826
$ws.Range("H132").Value = 1339.1538
$ws.Range("I132").Value = 1122.5
$ws.Range("K132").Value = 10102.5
$ws.Range("M132").Value = -7572.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 766.6667
$ws.Range("I23").Value = 766.6667
$ws.Range("K23").Value = 766.6667
$ws.Range("M23").Value = -543.6667
$ws.Range("H97").Value = 887.7368
$ws.Range("I97").Value = 789.86206
$ws.Range("J97").Value = 1203.1111
$ws.Range("K97").Value = 789.86206
$ws.Range("L97").Value = 1203.1111
$ws.Range("M97").Value = -293.86206
$ws.Range("N97").Value = -2195.1111
$ws.Range("H102").Value = 33468.375
$ws.Range("I102").Value = 1739.64
$ws.Range("K102").Value = 1739.64
$ws.Range("M102").Value = -117.6400000000001
$ws.Range("H113").Value = 5395
$ws.Range("I113").Value = 6992.3335
$ws.Range("J113").Value = 2999
$ws.Range("K113").Value = 6992.3335
$ws.Range("L113").Value = 2999
$ws.Range("M113").Value = -4822.3335
$ws.Range("N113").Value = -7339
$ws.Range("H122").Value = 3992
$ws.Range("I122").Value = 4117.909
$ws.Range("K122").Value = 12353.727
$ws.Range("M122").Value = -9903.726999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7282.1714
$ws.Range("I40").Value = 7378.706
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 7378.706
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -7242.706
$ws.Range("N40").Value = -4272
$ws.Range("H61").Value = 5725.52
$ws.Range("J61").Value = 6818.6665
$ws.Range("L61").Value = 6818.6665
$ws.Range("N61").Value = -7222.6665
$ws.Range("H102").Value = 94499.5
$ws.Range("J102").Value = 94499.5
$ws.Range("L102").Value = 94499.5
$ws.Range("N102").Value = -100989.5
$ws.Range("H113").Value = 5725.52
$ws.Range("J113").Value = 6818.6665
$ws.Range("L113").Value = 6818.6665
$ws.Range("N113").Value = -11158.6665
$ws.Range("H132").Value = 3018.55
$ws.Range("I132").Value = 2780.2144
$ws.Range("J132").Value = 3574.6667
$ws.Range("K132").Value = 8340.643199999999
$ws.Range("L132").Value = 10724.0001
$ws.Range("M132").Value = -5810.643199999999
$ws.Range("N132").Value = -15784.0001
$ws.Range("H136").Value = 5441.4644
$ws.Range("I136").Value = 5457.815
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 16373.445
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -13823.445
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6410.923
$ws.Range("I136").Value = 6650.4585
$ws.Range("K136").Value = 19951.3755
$ws.Range("M136").Value = -17401.3755
